$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C contains the "Förändrad" (Changed) date, stored as the Excel
# serial date number 46060 for every data row (2-79). Bump it by one day
# (46061) for each row, matching the automatic update in the diff.
for ($row = 2; $row -le 79; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
